{"js": "// Replace the three-digit-by-one-digit multiplication expressions in the\n// table cells with the new values from the commit.\nconst replacements = [\n    [\"441\u00d74=1764\", \"126\u00d77=882\"],\n    [\"380\u00d72=760\", \"467\u00d73=1401\"],\n    [\"818\u00d75=4090\", \"806\u00d72=1612\"],\n    [\"300\u00d76=1800\", \"147\u00d79=1323\"],\n    [\"984\u00d79=8856\", \"346\u00d73=1038\"],\n    [\"442\u00d77=3094\", \"756\u00d77=5292\"],\n    [\"472\u00d75=2360\", \"516\u00d73=1548\"],\n    [\"217\u00d75=1085\", \"158\u00d73=474\"],\n    [\"879\u00d75=4395\", \"779\u00d79=7011\"],\n    [\"969\u00d75=4845\", \"369\u00d73=1107\"],\n    [\"349\u00d77=2443\", \"992\u00d79=8928\"],\n    [\"842\u00d75=4210\", \"788\u00d79=7092\"],\n    [\"176\u00d77=1232\", \"526\u00d76=3156\"],\n    [\"846\u00d72=1692\", \"878\u00d79=7902\"],\n    [\"823\u00d77=5761\", \"847\u00d73=2541\"],\n    [\"509\u00d76=3054\", \"782\u00d75=3910\"],\n    [\"546\u00d76=3276\", \"850\u00d76=5100\"],\n    [\"179\u00d79=1611\", \"223\u00d72=446\"],\n    [\"275\u00d75=1375\", \"252\u00d77=1764\"],\n    [\"780\u00d78=6240\", \"485\u00d78=3880\"],\n    [\"470\u00d78=3760\", \"218\u00d79=1962\"],\n    [\"811\u00d77=5677\", \"823\u00d79=7407\"],\n    [\"240\u00d78=1920\", \"508\u00d78=4064\"],\n    [\"383\u00d75=1915\", \"848\u00d72=1696\"],\n    [\"947\u00d73=2841\", \"844\u00d77=5908\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n    const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n    results.load(\"items\");\n    await context.sync();\n\n    for (let i = 0; i < results.items.length; i++) {\n        results.items[i].insertText(newText, Word.InsertLocation.replace);\n    }\n    await context.sync();\n}\n", "ps1": "# Replace the three-digit-by-one-digit multiplication expressions in the\n# table cells with the new values from the commit.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"441\u00d74=1764\", \"126\u00d77=882\"),\n    @(\"380\u00d72=760\", \"467\u00d73=1401\"),\n    @(\"818\u00d75=4090\", \"806\u00d72=1612\"),\n    @(\"300\u00d76=1800\", \"147\u00d79=1323\"),\n    @(\"984\u00d79=8856\", \"346\u00d73=1038\"),\n    @(\"442\u00d77=3094\", \"756\u00d77=5292\"),\n    @(\"472\u00d75=2360\", \"516\u00d73=1548\"),\n    @(\"217\u00d75=1085\", \"158\u00d73=474\"),\n    @(\"879\u00d75=4395\", \"779\u00d79=7011\"),\n    @(\"969\u00d75=4845\", \"369\u00d73=1107\"),\n    @(\"349\u00d77=2443\", \"992\u00d79=8928\"),\n    @(\"842\u00d75=4210\", \"788\u00d79=7092\"),\n    @(\"176\u00d77=1232\", \"526\u00d76=3156\"),\n    @(\"846\u00d72=1692\", \"878\u00d79=7902\"),\n    @(\"823\u00d77=5761\", \"847\u00d73=2541\"),\n    @(\"509\u00d76=3054\", \"782\u00d75=3910\"),\n    @(\"546\u00d76=3276\", \"850\u00d76=5100\"),\n    @(\"179\u00d79=1611\", \"223\u00d72=446\"),\n    @(\"275\u00d75=1375\", \"252\u00d77=1764\"),\n    @(\"780\u00d78=6240\", \"485\u00d78=3880\"),\n    @(\"470\u00d78=3760\", \"218\u00d79=1962\"),\n    @(\"811\u00d77=5677\", \"823\u00d79=7407\"),\n    @(\"240\u00d78=1920\", \"508\u00d78=4064\"),\n    @(\"383\u00d75=1915\", \"848\u00d72=1696\"),\n    @(\"947\u00d73=2841\", \"844\u00d77=5908\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
